$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "designer"
$ws.Range("B4").Value = 12345
$ws.Range("C4").Value = 12345
$ws.Range("D4").Value = "Magasin Ikea"
$ws.Range("E4").Value = "Test with numeric values"
$ws.Range("G4").Value = "blue"
$ws.Range("H4").Value = 12345
$ws.Range("I4").Value = "blue,green"

$ws.Range("E5").Select()
